$wb = $excel.ActiveWorkbook

# --- Worksheet references ---
$ws1 = $wb.Worksheets.Item(1)   # Weekly Quantity
$ws2 = $wb.Worksheets.Item(2)   # Monthly Trend

# --- Rename header labels on the existing sheets ---
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the existing sheets ---
$lastIndex = $wb.Worksheets.Count
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$ws3.Name = "PO Forecast"

# --- Header row ---
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Reuse the same header formatting (bold, centered, bordered) used on the
# other sheets, and the same date number format for column A.
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:B1").PasteSpecial(-4122)
$ws3.Range("C1:D1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws3.Range("A2:A12").PasteSpecial(-4122)

# --- Forecast data rows ---
$ws3.Range("A2").Value = 44948.99999999999
$ws3.Range("B2").Value = 5
$ws3.Range("C2").Value = 3.220592490815414
$ws3.Range("D2").Value = 6.728880059427539

$ws3.Range("A3").Value = 44955.99999999999
$ws3.Range("B3").Value = 8
$ws3.Range("C3").Value = 6.323193263148783
$ws3.Range("D3").Value = 9.740928524968536

$ws3.Range("A4").Value = 44962.99999999999
$ws3.Range("B4").Value = 11
$ws3.Range("C4").Value = 9.192093891122026
$ws3.Range("D4").Value = 12.80641442163465

$ws3.Range("A5").Value = 44969.99999999999
$ws3.Range("B5").Value = 14
$ws3.Range("C5").Value = 12.19524784476156
$ws3.Range("D5").Value = 15.92724707188889

$ws3.Range("A6").Value = 44976.99999999999
$ws3.Range("B6").Value = 17
$ws3.Range("C6").Value = 15.18020060654142
$ws3.Range("D6").Value = 18.76430587235053

$ws3.Range("A7").Value = 44983.99999999999
$ws3.Range("B7").Value = 20
$ws3.Range("C7").Value = 18.24451118627795
$ws3.Range("D7").Value = 21.78397619660839

$ws3.Range("A8").Value = 44990.99999999999
$ws3.Range("B8").Value = 23
$ws3.Range("C8").Value = 21.17852693386302
$ws3.Range("D8").Value = 24.91967454449661

$ws3.Range("A9").Value = 44997.99999999999
$ws3.Range("B9").Value = 26
$ws3.Range("C9").Value = 24.3829193615796
$ws3.Range("D9").Value = 27.77553422730972

$ws3.Range("A10").Value = 45004.99999999999
$ws3.Range("B10").Value = 29
$ws3.Range("C10").Value = 27.25983500299118
$ws3.Range("D10").Value = 30.7539376006356

$ws3.Range("A11").Value = 45011.99999999999
$ws3.Range("B11").Value = 32
$ws3.Range("C11").Value = 30.2189888564153
$ws3.Range("D11").Value = 33.74057111868668

$ws3.Range("A12").Value = 45018.99999999999
$ws3.Range("B12").Value = 35
$ws3.Range("C12").Value = 33.32277240844854
$ws3.Range("D12").Value = 36.77458879917329

$ws3.Range("A1").Select()

# Restore the originally active sheet/selection.
$ws1.Activate()
$ws1.Range("A1").Select()
